# Creacion de Excel BD
# Adds earthquake records (rows 26-29) to the "Hoja1" sheet and updates
# the "Descripcion Detallada" note on row 2 (J2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update J2 ("Descripcion Detallada" for the first record) ---
$ws.Range("J2").Value = "Sucedió"

# --- Row 26 ---
$ws.Range("C26").NumberFormat = "@"
$ws.Range("F26:H26").NumberFormat = "@"
$ws.Range("A26").Value = "24/01/2016"
$ws.Range("B26").Value = "05:02:45"
$ws.Range("C26").Value = "4.55"
$ws.Range("D26").Value = "TECTONICO_SUBDUCCION"
$ws.Range("E26").Value = "detalle"
$ws.Range("F26").Value = "89.0"
$ws.Range("G26").Value = "9.7"
$ws.Range("H26").Value = "6.0"
$ws.Range("I26").Value = "ALAJUELA"
$ws.Range("J26").Value = "Descripcion"

# --- Row 27 ---
$ws.Range("C27").NumberFormat = "@"
$ws.Range("F27:H27").NumberFormat = "@"
$ws.Range("A27").Value = "24/12/2016"
$ws.Range("B27").Value = "05:02:45"
$ws.Range("C27").Value = "4.55"
$ws.Range("D27").Value = "TECTONICO_SUBDUCCION"
$ws.Range("E27").Value = "detalle"
$ws.Range("F27").Value = "89.8"
$ws.Range("G27").Value = "5.0"
$ws.Range("H27").Value = "6.0"
$ws.Range("I27").Value = "ALAJUELA"
$ws.Range("J27").Value = "Descripcion"

# --- Row 28 ---
$ws.Range("C28").NumberFormat = "@"
$ws.Range("F28:H28").NumberFormat = "@"
$ws.Range("A28").Value = "24/07/2017"
$ws.Range("B28").Value = "23:07:46"
$ws.Range("C28").Value = "8.7"
$ws.Range("D28").Value = "DEFORMACION_INTERNA"
$ws.Range("E28").Value = "dt"
$ws.Range("F28").Value = "5.0"
$ws.Range("G28").Value = "4.9"
$ws.Range("H28").Value = "67.2"
$ws.Range("I28").Value = "GUANACASTE"
$ws.Range("J28").Value = "Descripcion"

# --- Row 29 ---
$ws.Range("C29").NumberFormat = "@"
$ws.Range("F29:H29").NumberFormat = "@"
$ws.Range("A29").Value = "24/01/2020"
$ws.Range("B29").Value = "20:02:45"
$ws.Range("C29").Value = "4.0"
$ws.Range("D29").Value = "TECTONICO_POR_FALLA_LOCAL"
$ws.Range("E29").Value = "detalle"
$ws.Range("F29").Value = "89.0"
$ws.Range("G29").Value = "9.7"
$ws.Range("H29").Value = "6.0"
$ws.Range("I29").Value = "SAN_JOSE"
$ws.Range("J29").Value = "Natalia"
